$wb = $excel.ActiveWorkbook

# --- Sheet "BookShelves": add 3 new bookshelf product rows ---
$wsBook = $wb.Worksheets.Item("BookShelves")

$wsBook.Range("A2").Value = "Rhodes Solid Wood Bookshelf In Mahogany Finish"
$wsBook.Range("B2").Value = "By Urban Ladder"
$wsBook.Range("C2").Value = "₹14,755"

$wsBook.Range("A3").Value = "Rhodes Solid Wood Bookshelf In Teak Finish"
$wsBook.Range("B3").Value = "By Urban Ladder"
$wsBook.Range("C3").Value = "₹14,755"

$wsBook.Range("A4").Value = "Theodore Engineered Wood Bookshelf In Rustic Walnut Finish"
$wsBook.Range("B4").Value = "By Urban Ladder"
$wsBook.Range("C4").Value = "₹12,814"

# --- Sheet "submenuItems": add the full category submenu list ---
$wsMenu = $wb.Worksheets.Item("submenuItems")

$wsMenu.Range("A2").Value = "Lounge Chairs"
$wsMenu.Range("B2").Value = "TV Units"

$wsMenu.Range("A3").Value = "Accent Chairs"
$wsMenu.Range("B3").Value = "Bookshelves"

$wsMenu.Range("A4").Value = "Recliners"
$wsMenu.Range("B4").Value = "Shoe Racks"

$wsMenu.Range("A5").Value = "Sofa Cum Bed"
$wsMenu.Range("B5").Value = "Prayer Units"

$wsMenu.Range("A6").Value = "UL Assured Picks"
$wsMenu.Range("B6").Value = "Showcases"

$wsMenu.Range("A7").Value = "Ottomans & Stools"
$wsMenu.Range("B7").Value = "Wall Shelves"

$wsMenu.Range("A8").Value = "Bean Bags"
$wsMenu.Range("B8").Value = "Entryway & Foyer"

$wsMenu.Range("A9").Value = "Benches"
$wsMenu.Range("B9").Value = "Room Divider"

$wsMenu.Range("A10").Value = "Bar Stools"
$wsMenu.Range("B10").Value = "Living Room Sets"

$wsMenu.Range("A11").Value = "Rocking Chairs"

$wsMenu.Range("A12").Value = "Gaming Chairs"

# --- Sheet "StudyChairs": blank out Brand Name / Price for row 2 ---
$wsChairs = $wb.Worksheets.Item("StudyChairs")
$wsChairs.Range("B2").Value = ""
$wsChairs.Range("C2").Value = ""
